# "Nut roast up date." — add a new "marmite" ingredient row to the
# Ingredients sheet, keeping the alphabetical ordering (it sits between
# "lime pickle" and "marjoram", i.e. at row 106), which pushes every
# following row down by one (183 -> 184 total data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 106 - this shifts rows 106..183 down to 107..184.
$ws.Rows.Item(106).Insert()

# Populate the new ingredient row (same shape as the other "Check" type
# ingredients that carry no nutritional data yet).
$ws.Cells.Item(106, 1).Value = "marmite"
$ws.Cells.Item(106, 2).Value = "Check"
$ws.Cells.Item(106, 3).Value = 0
$ws.Cells.Item(106, 4).Value = 0
$ws.Cells.Item(106, 5).Value = 0
$ws.Cells.Item(106, 6).Value = 0
$ws.Cells.Item(106, 7).Value = 0

# The sheet's used range/filter range grew by the one inserted row.
$wb.Names.Item("Ingredients!_FilterDatabase").RefersTo = "=Ingredients!`$A`$2:`$G`$184"

# Reflect where the author was working when they saved: scrolled down to
# the newly added row and with its data cells selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 86
[void]$ws.Range("C106:G106").Select()
